$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'46.929.64"
$ws.Range("E2").Value = "'  +5.62%  "
$ws.Range("D3").Value = "'2.337.66"
$ws.Range("E3").Value = "'  +5.04%  "
$ws.Range("E4").Value = "'  -0.77%  "
$ws.Range("D5").Value = "'305.66"
$ws.Range("E5").Value = "'  +0.85%  "
$ws.Range("D6").Value = "'97.58"
$ws.Range("E6").Value = "'  +8.27%  "
$ws.Range("E7").Value = "'  +3.65%  "
$ws.Range("E8").Value = "'  -0.68%  "
$ws.Range("E9").Value = "'  +8.39%  "
$ws.Range("D10").Value = "'35.97"
$ws.Range("E10").Value = "'  +6.50%  "
$ws.Range("D11").Value = "'0.0814"
$ws.Range("E11").Value = "'  +3.93%  "
$ws.Range("D12").Value = "'7.47"
$ws.Range("E12").Value = "'  +7.36%  "
$ws.Range("E13").Value = "'  -0.10%  "
$ws.Range("D14").Value = "'2.691.59"
$ws.Range("E14").Value = "'  +4.84%  "
$ws.Range("D15").Value = "'2.333.85"
$ws.Range("E15").Value = "'  +0.54%  "
$ws.Range("D16").Value = "'14.19"
$ws.Range("E16").Value = "'  +7.45%  "
$ws.Range("D17").Value = "'0.836"
$ws.Range("E17").Value = "'  +3.48%  "
$ws.Range("D18").Value = "'46.780.92"
$ws.Range("E18").Value = "'  +5.81%  "
$ws.Range("D19").Value = "'13.68"
$ws.Range("E19").Value = "'  +19.73%  "
$ws.Range("E20").Value = "'  +4.54%  "
$ws.Range("D21").Value = "'6.22"
$ws.Range("E21").Value = "'  +2.97%  "
$ws.Range("D22").Value = "'67.89"
$ws.Range("E22").Value = "'  +5.30%  "
$ws.Range("D23").Value = "'251.44"
$ws.Range("E23").Value = "'  +7.60%  "
$ws.Range("E24").Value = "'  +4.20%  "
$ws.Range("E25").Value = "'  +4.33%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "'  -0.58%  "
$ws.Range("D27").Value = "'42.63"
$ws.Range("E27").Value = "'  +17.59%  "
$ws.Range("D28").Value = "'2.27"
$ws.Range("E28").Value = "'  +0.10%  "
$ws.Range("D29").Value = "'9.91"
$ws.Range("E29").Value = "'  +4.35%  "
$ws.Range("D30").Value = "'20.34"
$ws.Range("E30").Value = "'  +3.89%  "
$ws.Range("D31").Value = "'5.83"
$ws.Range("E31").Value = "'  +4.09%  "
$ws.Range("D32").Value = "'0.0817"
$ws.Range("E32").Value = "'  +7.79%  "
$ws.Range("D33").Value = "'146.47"
$ws.Range("E33").Value = "'  +0.12%  "
$ws.Range("E34").Value = "'  +0.01%  "
$ws.Range("B35").Value = "'Kaspa"
$ws.Range("C35").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.115"
$ws.Range("E35").Value = "'  +7.76%  "
$ws.Range("B36").Value = "'LidoDAOToken"
$ws.Range("C36").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "'3.15"
$ws.Range("E36").Value = "'  +4.57%  "
$ws.Range("E37").Value = "'  +3.45%  "
$ws.Range("E38").Value = "'  +1.81%  "
$ws.Range("D39").Value = "'4.00"
$ws.Range("E39").Value = "'  +9.37%  "
$ws.Range("D40").Value = "'0.0312"
$ws.Range("E40").Value = "'  +8.68%  "
$ws.Range("D41").Value = "'3.40"
$ws.Range("E41").Value = "'  +5.21%  "
$ws.Range("D42").Value = "'14.09"
$ws.Range("E42").Value = "'  -3.69%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "'  -0.86%  "
$ws.Range("D44").Value = "'1.96"
$ws.Range("E44").Value = "'  +13.49%  "
$ws.Range("D45").Value = "'1.808.29"
$ws.Range("E45").Value = "'  +1.41%  "
$ws.Range("D46").Value = "'92.17"
$ws.Range("E46").Value = "'  +15.72%  "
$ws.Range("D47").Value = "'74.81"
$ws.Range("E47").Value = "'  +11.43%  "
$ws.Range("D48").Value = "'0.196"
$ws.Range("E48").Value = "'  +7.35%  "
$ws.Range("D49").Value = "'98.98"
$ws.Range("E49").Value = "'  +3.71%  "
$ws.Range("D50").Value = "'55.37"
$ws.Range("E50").Value = "'  +5.72%  "
$ws.Range("E51").Value = "'  +5.59%  "
